# Update the self-introduction slide (slide 2) profile bullet lines
# to the new GAS-appropriate speaker profile items.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$s.Shapes.Item(17).TextFrame.TextRange.Text = "▶  非エンジニア × AI × GASで業務改善！"
$s.Shapes.Item(18).TextFrame.TextRange.Text = "▶  Notion公式アンバサダー"
$s.Shapes.Item(19).TextFrame.TextRange.Text = "▶  Notionさいたま主宰 / DATASaber"
$s.Shapes.Item(20).TextFrame.TextRange.Text = "▶  X: @keitaro_aigc"
